$d = $word.ActiveDocument

function New-WordPkgXml($bodyInner) {
    return '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
           '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
           '<pkg:xmlData>' +
           '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
           '<w:body>' + $bodyInner + '</w:body>' +
           '</w:document>' +
           '</pkg:xmlData></pkg:part></pkg:package>'
}

# Content (as raw w:p OOXML) for the five new paragraphs that need to be
# inserted right before the final ("_GoBack" bookmarked) paragraph.
$newParas = @(
    '<w:p><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>2</w:t></w:r><w:r><w:t xml:space="preserve">. </w:t></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>优化了热插拔，这包括p</w:t></w:r><w:r><w:t>lugin manager, system manager</w:t></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>及v</w:t></w:r><w:r><w:t>iew editor</w:t></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>。</w:t></w:r></w:p>',
    '<w:p/>',
    '<w:p><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>2020.09.28</w:t></w:r></w:p>',
    '<w:p><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>1.</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>优化了各个插件U</w:t></w:r><w:r><w:t>I</w:t></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>的自适应。</w:t></w:r></w:p>',
    '<w:p><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>2.</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>为控制台版本增加了获取系统信息命令。</w:t></w:r></w:p>'
)

# Step 1: insert five empty paragraphs immediately before the very last
# paragraph of the document (the one that carries the "_GoBack" bookmark),
# then stamp each with its final OOXML content via InsertXML.
foreach ($paraXml in $newParas) {
    $last = $d.Paragraphs.Last
    $insPos = $last.Range.Duplicate
    $insPos.Collapse(1)
    [void]$insPos.InsertParagraphBefore()

    $target = $d.Paragraphs.Item($d.Paragraphs.Count - 1)
    [void]$target.Range.InsertXML((New-WordPkgXml $paraXml))
}

# Step 2: rewrite the run content of the final paragraph (keep the
# paragraph's pPr and the "_GoBack" bookmark untouched/at the end).
$last = $d.Paragraphs.Last
$lastRange = $last.Range
$start = $lastRange.Start
$end = $lastRange.End

# Remove the existing run text (everything up to, but excluding, the
# paragraph mark) so the bookmark, which is anchored right before the
# paragraph mark, collapses down to the (now empty) insertion point.
[void]$d.Range($start, $end - 1).Delete()

$newRunsXml = '<w:p><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>3.</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>修复了部分b</w:t></w:r><w:r><w:t>ug</w:t></w:r></w:p>'
[void]$d.Range($start, $start).InsertXML((New-WordPkgXml $newRunsXml))
